$d = $word.ActiveDocument

# --- Step 1: merge the 3 runs of the first (and only) paragraph into a single
# run with the complete sentence. Using Find/Execute with matching Find/Replace
# text (an effective no-op textually) rewrites the whole range as one run,
# which conveniently drops the <w:proofErr> spell-check markers that bracketed
# "Git" while leaving everything else (incl. the _GoBack bookmark) untouched.
$fullSentence = "This is my test of a binary file being tracked in a Git repository."
$d.Content.Find.Execute($fullSentence, $false, $false, $false, $false, $false, `
    $true, 1, $false, $fullSentence, 2) | Out-Null

# --- Step 2: relocate the _GoBack bookmark from the end of paragraph 1 to the
# end of the (new) paragraph 2, after the new run of text.
$bm = $d.Bookmarks("_GoBack")

# Insert the paragraph break right where the bookmark currently sits; the
# bookmark itself remains anchored at the end of paragraph 1.
$breakPoint = $d.Range($bm.Start, $bm.Start)
$breakPoint.InsertParagraphAfter()

# Drop the now out-of-place bookmark; we'll recreate it in the right spot.
$d.Bookmarks("_GoBack").Delete()

# --- Step 3: populate the new (second) paragraph with the added sentence,
# plus one temporary placeholder character at the end.
$para2 = $d.Paragraphs(2).Range
$para2.InsertBefore("I HAVE CHANGED THE FILE!!X")

# --- Step 4: re-create the _GoBack bookmark. A *collapsed* (zero-width) range
# placed exactly at a paragraph end trips a quirk in this host's bookmark
# engine, so instead we wrap the temporary trailing "X" placeholder (a proper,
# non-collapsed range) and then delete that character through the bookmark's
# own Range, which cleanly collapses the bookmark to zero width in place -
# landing right after the real text and before the paragraph mark.
$para2Fresh = $d.Paragraphs(2).Range
$null = $para2Fresh.MoveEnd(1, -1)
$placeholder = $d.Range($para2Fresh.End - 1, $para2Fresh.End)
$d.Bookmarks.Add("_GoBack", $placeholder)

$bmRange = $d.Bookmarks("_GoBack").Range
$bmRange.Text = ""
